# Generate Report for Handback
# Updates the handoff/handback timestamps for the
# ee154e4d-69c6-475a-83f5-9caf1fb6f52e.md file row across the
# zh-cn / de-de status sheets, and refreshes the corresponding
# "Latest HO Xliff Generate Date" on the Overview sheet.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-20 06:53:25"
$wsZhCn.Range("K3").Value = "2016-08-20 06:53:41"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-20 06:53:29"
$wsDeDe.Range("K3").Value = "2016-08-20 06:53:47"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-20 06:53:29"
